# Apply the corrections described by the commit:
#  - Fix a truncated/garbled closing bracket in three duplicated text blocks
#    (R21 - PARENTHESES rule text): "...is a design decision.} " ->
#    "...is a design decision. " (drop the stray "}").
#  - Because that text shrank by one character, every downstream character
#    offset recorded in columns B/C/F/G (rows 22-43) shifts down by one.
#  - Remove the leftover cell formatting on L3:L4 that isn't used anywhere
#    else in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the garbled text in H22, J22, L22 -------------------------------
$textCells = @("H22", "J22", "L22")
foreach ($addr in $textCells) {
    $cell = $ws.Range($addr)
    $val = $cell.Value2
    $cell.Value = $val.Replace("design decision.}", "design decision.")
}

# Re-run autofit on the affected row so no stray explicit row height sticks
# around as a side effect of the in-place text edit above.
$ws.Rows("22:22").AutoFit()

# --- 2. Shift the downstream character offsets by -1 ------------------------
$ws.Range("G22").Value = $ws.Range("G22").Value2 - 1

for ($r = 23; $r -le 43; $r++) {
    foreach ($col in @("B", "C", "F", "G")) {
        $addr = "$col$r"
        $cur = $ws.Range($addr).Value2
        $ws.Range($addr).Value = $cur - 1
    }
}

# --- 3. Drop the unused/leftover style applied to L3:L4 ---------------------
$ws.Range("L3:L4").ClearFormats()
